# Generate Report for Handoff
# Inserts a new "45a260cc-ee79-46db-8053-8fce3241eb47.md" row ahead of the
# existing "cffef19d-ae73-412d-b190-eb723fa21bd1.md" row on every sheet
# (Overview, zh-cn, de-de), pushing the old row from row 2 down to row 3.

$wb = $excel.ActiveWorkbook

$commitHash = "ae3f0e64a851b9e2d85277f3203cf576cab72ebc"
$oldFile = "cffef19d-ae73-412d-b190-eb723fa21bd1.md"
$newFile = "45a260cc-ee79-46db-8053-8fce3241eb47.md"

function Set-GithubHyperlink($ws, $cellAddr, $fileName, $displayPrefix) {
    $url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileName"
    $display = "$displayPrefix$fileName"
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A..G = File Name, Path And Name, Extension,
# Publish URL, zh-cn, de-de, Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Rows.Item(2).Insert()

# drop the stale hyperlink left behind on B2 by the row insert, then
# recreate both hyperlinks (new file keeps rId2, old file gets rId3)
$wsOv.Range("B2").Hyperlinks.Delete()

$wsOv.Range("A2").Value = $newFile
$wsOv.Range("B2").Value = "e2e\$newFile"
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("E2").Value = "Ready for handoff"
$wsOv.Range("F2").Value = "Ready for handoff"
$wsOv.Range("G2").Value = "2016-08-19 00:38:47"
$wsOv.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-GithubHyperlink $wsOv "B2" $newFile "e2e\"
Set-GithubHyperlink $wsOv "B3" $oldFile "e2e\"

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A..P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(2).Insert()

$wsZh.Range("A2").Hyperlinks.Delete()

$wsZh.Range("A2").Value = $newFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "45a260cc-ee79-46db-8053-8fce3241eb47.a61ab196186629d752af02e2934727b770208b80.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-19 00:38:42"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M2").Value = "True"
$wsZh.Range("O2").Value = "False"

Set-GithubHyperlink $wsZh "A2" $newFile ""
Set-GithubHyperlink $wsZh "A3" $oldFile ""

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A..P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(2).Insert()

$wsDe.Range("A2").Hyperlinks.Delete()

$wsDe.Range("A2").Value = $newFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "45a260cc-ee79-46db-8053-8fce3241eb47.a61ab196186629d752af02e2934727b770208b80.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-19 00:38:47"
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M2").Value = "True"
$wsDe.Range("O2").Value = "False"

Set-GithubHyperlink $wsDe "A2" $newFile ""
Set-GithubHyperlink $wsDe "A3" $oldFile ""

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
